# Commit: "Add files via upload" -- update the confirmation checkbox
# field definition on the "Tabelle1" sheet (row 33, column D) to add
# two extra checkbox fields ("sonstwas" and "irgendwas"), matching a
# fresh upload of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Core content change: D33 field-definition string gains two more
# checkbox entries.
$ws.Range("D33").Value = "Richtig und Vollständig:checkbox*;sonstwas:checkbox;irgendwas:checkbox"

# The wider text means column D is now sized to fit it (bestFit column).
$ws.Columns.Item(4).ColumnWidth = 218.666666666667

# Reflect the author's new selection on that same cell.
$ws.Range("D33").Select()
